$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D4").Value = "2016-01-18 06:31:24"
$zhcn.Range("G4").Value = "2016-01-18 06:32:07"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D4").Value = "2016-01-18 06:31:34"
$dede.Range("G4").Value = "2016-01-18 06:32:23"
